$d = $word.ActiveDocument

# --- 1) Merge "1" + "ª Entrevista " into a single run "1ª Entrevista"
#        (drops the trailing space) and underline the whole heading. ---
$null = $d.Content.Find.Execute("1ª Entrevista ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "1ª Entrevista", 2)

$headingRange = $d.Paragraphs(1).Range
$headingRange.Font.Underline = 1

# --- 2) Move the "_GoBack" last-edit bookmark onto this heading. ---
# The heading text is now exactly "1ª Entrevista" (13 characters), so the
# spot right after the run sits at the very end of the paragraph's text
# (right before the paragraph mark). Dropping a bookmark directly there
# snaps it around the whole paragraph, so instead: temporarily extend the
# paragraph by one placeholder character, plant the bookmark right before
# that placeholder, then delete the placeholder again - leaving the
# bookmark correctly seated right after the run and before the paragraph
# mark. The end position is cached in a plain number up front so it
# does not drift as the document is mutated.
$endPos = $headingRange.End - 1

$placeholderSpot = $d.Range($endPos, $endPos)
$placeholderSpot.InsertAfter("X")

$bookmarkSpot = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$placeholder = $d.Range($endPos, $endPos + 1)
$placeholder.Delete()
